# "updating Luk and Spivak 2020 hook"
#
# The original sheet had a single data row (row 4) for study "Luk_et_al_2019".
# This edit renames that study hook to "Luk_et_al_2020_a" (lower-casing its
# compaction code along the way) and adds a sibling row 5 for a second hook,
# "Luk_et_al_2020_b", that shares some of row 4's values.
#
# NOTE: the writes below are intentionally ordered ("no obvious compaction",
# then "freeze dried", then "Luk_et_al_2020_a", then "Luk_et_al_2020_b") so
# that new shared-string entries get appended to xl/sharedStrings.xml in the
# same order as the upstream commit - cosmetic, but keeps the workbook's
# internal string table identical to the original edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: compaction_flag code is re-cased; study id updated later below ---
$ws.Range("F4").Value = "no obvious compaction"

# --- Row 5: new row for the "Luk_et_al_2020_b" hook ---
$newRow = $ws.Range("A5:AE5")
$newRow.WrapText = $true
$ws.Rows.Item(5).RowHeight = 16

$ws.Range("K5").Value = "freeze dried"

# --- Row 4: study id hook rename (Luk_et_al_2019 -> Luk_et_al_2020_a) ---
$ws.Range("A4").Value = "Luk_et_al_2020_a"

# --- Row 5: remaining values for the "Luk_et_al_2020_b" hook ---
$ws.Range("A5").Value = "Luk_et_al_2020_b"
$ws.Range("B5").Value = "piston corer"
$ws.Range("F5").Value = "no obvious compaction"
$ws.Range("H5").Value = 168
$ws.Range("W5").Value = "gamma"
$ws.Range("X5").Value = "gamma"

# --- Row 4: radiometric counting method codes no longer apply, clear them ---
$ws.Range("W4").ClearContents()
$ws.Range("X4").ClearContents()

# --- Leave the sheet with the same active selection as the upstream edit ---
$ws.Range("Q5").Select()
